$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Update the three body-text paragraphs in "TextBox 12" (shape 3) ---
$sh = $s.Shapes.Item(3)
$tr = $sh.TextFrame.TextRange

function Replace-TextRange($textRange, [string]$oldText, [string]$newText) {
    $full = $textRange.Text
    $idx = $full.IndexOf($oldText)
    if ($idx -lt 0) {
        throw "Could not find target text: $oldText"
    }
    $sub = $textRange.Characters($idx + 1, $oldText.Length)
    $sub.Text = $newText
}

Replace-TextRange $tr 'For the rover''s GPS Navigation functions, we are using an algorithm that determines the shortest path between two given GPS coordinates. The GPS will also keep updating new best routes per request from obstacle avoidance and unstack functions. Which means that the GPS function has to work flawlessly with both of these two functions to ensure the rover''s safety and efficiency. How rover behaves during its driving is also critical, the GPS function will check if the rover is off-course by preset time interval and give route compensation if needed.  ' 'For the rover’s GPS Navigation functions, we are using an algorithm that determines the shortest path between two given GPS coordinates. The GPS will also keep updating the new best route per request from the obstacle avoidance and unstuck from obstacles modules. This means that the GPS function has to work flawlessly with both of these modules to ensure the rover’s safety and efficiency. How the rover behaves during it’s driving is also critical, so the GPS function will check if the rover is off-course every few seconds and give route compensation if needed.  '
Replace-TextRange $tr 'The obstacle avoidance system ensures that our rover is not impeded on its way to the destination. Taking in filtered images from the obstacle detection software, this system does edge detection on the image find objects in the rovers path, and then decides how to best get around the object. This is done by treating the filtered black and white image as a matrix of pixels, and summing the number of edges to the left, right or in front of the rover and adjusting the direction of the rover to travel where the fewest edges are found.' 'The obstacle avoidance system ensures that our rover is not impeded on its way to the destination. Taking in filtered images from the obstacle detection software, this system does edge detection on the image to find objects in the rovers path, and then decides how to best get around the object. This is done by treating the filtered black and white image as a matrix of pixels, and summing the number of edges to the left, right or in front of the rover and adjusting the direction of the rover to travel where the fewest edges are found.'
Replace-TextRange $tr 'Once the rover get''s within the GPS'' error range of the finish coordinates, we have to search for the finish pole. This algorithm works by first searching for the finish pole by rotating in place and taking pictures, then aligning the rover in the direction of the finish, and moving forward, making periodic course corrections along the way. ' 'Once the rover get’s within the GPS’ error range of the finish coordinates, we have to search for the finish pole. This algorithm works by first searching for the finish pole by rotating in place and taking pictures. These pictures are used to detect a traffic cone by our imaging system. Once the cone is detected, the rover is oriented in the direction of the cone, and moves forward, making periodic course corrections along the way. '

# --- Resize "TextBox 12" (shape 3) to its new height ---
$sh.Left = 981.1364566929134
$sh.Top = 335.971811023622
$sh.Width = 726.1791338582677
$sh.Height = 1799.684409448819

# --- Reposition "Picture 11" (shape 13) ---
$sh13 = $s.Shapes.Item(13)
$sh13.Left = 1073.8933070866142
$sh13.Top = 2135.656220472441

# --- Update the date placeholder fields on every slide layout + the slide master + notes master ---
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $layout = $p.SlideMaster.CustomLayouts.Item($i)
    for ($j = 1; $j -le $layout.Shapes.Count; $j++) {
        $lsh = $layout.Shapes.Item($j)
        if ($lsh.HasTextFrame -and $lsh.TextFrame.HasText) {
            $ltr = $lsh.TextFrame.TextRange
            if ($ltr.Text -eq '4/15/2017') {
                $ltr.Text = '4/17/2017'
            }
        }
    }
}

for ($j = 1; $j -le $p.SlideMaster.Shapes.Count; $j++) {
    $msh = $p.SlideMaster.Shapes.Item($j)
    if ($msh.HasTextFrame -and $msh.TextFrame.HasText) {
        $mtr = $msh.TextFrame.TextRange
        if ($mtr.Text -eq '4/15/2017') {
            $mtr.Text = '4/17/2017'
        }
    }
}

$nm = $p.NotesMaster
for ($j = 1; $j -le $nm.Shapes.Count; $j++) {
    $nsh = $nm.Shapes.Item($j)
    if ($nsh.HasTextFrame -and $nsh.TextFrame.HasText) {
        $ntr = $nsh.TextFrame.TextRange
        if ($ntr.Text -eq '4/15/2017') {
            $ntr.Text = '4/17/2017'
        }
    }
}
